$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values in column F (dSF) for specific rows, per repulled data.
$updates = @{
    5  = 3
    6  = -3
    8  = -5
    14 = -2
    15 = 4
    18 = -3
    20 = 2
    24 = -2
    29 = -4
    33 = -1
    34 = -1
    42 = -1
    43 = 2
    44 = 4
    45 = -3
    47 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
